$wb = $excel.ActiveWorkbook

# --- "cases" sheet: revert management (F) and rotation (G) columns for rows 2-6 ---
$ws4 = $wb.Worksheets.Item("cases")
for ($r = 2; $r -le 6; $r++) {
    $ws4.Cells.Item($r, 6).Value = '"WHEAT.Ble_Dur_1"'
    $ws4.Cells.Item($r, 7).Value = '"ROTATION_BLE_IRRIGUE"'
}

# --- restore previous selections on each sheet ---
$ws1 = $wb.Worksheets.Item("generalOptions")
$ws1.Range("B5").Select()

$ws2 = $wb.Worksheets.Item("testble mais poischiche")
$ws2.Range("B5").Select()

$ws3 = $wb.Worksheets.Item("caseswithoutmaize")
$ws3.Range("A3").Select()

$ws4.Range("H9").Select()
$ws4.Activate()
